# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 8636
$wsExhibit.Range("F3").Value = 89
$wsExhibit.Range("F5").Value = 93
$wsExhibit.Range("F6").Value = 1378
$wsExhibit.Range("F7").Value = 1370
$wsExhibit.Range("F8").Value = 231
$wsExhibit.Range("F9").Value = 34
$wsExhibit.Range("F10").Value = 253
$wsExhibit.Range("F11").Value = 76

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 8637
$wsAll.Range("F3").Value = 89
$wsAll.Range("F5").Value = 93
$wsAll.Range("F6").Value = 1378
$wsAll.Range("F7").Value = 1370
$wsAll.Range("F8").Value = 231
$wsAll.Range("F10").Value = 34
$wsAll.Range("F11").Value = 253
$wsAll.Range("F12").Value = 76
